$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap LEG1_DIRECTION (AL2) and LEG2_DIRECTION (AM2) values
$ws.Range("AL2").Value = "R"
$ws.Range("AM2").Value = "P"

# Reflect the user's final selection/scroll position after making the edit
$ws.Range("AM2").Select() | Out-Null
